# "Fixed BOM and PnP files"
#
# The LCSC part number for the FT232RL USB IC (row 10, "LCSC Part #" column,
# cell E10 on the BOM sheet) was wrong. Correct it from the old/defunct
# part C490691 to the correct part C8690 (FTDI FT232RL-REEL), and point the
# cell's hyperlink at the matching LCSC product page.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$cell = $ws.Range("E10")

# Update the visible BOM value.
$cell.Value = "C8690"

# Update the hyperlink target to the new LCSC product page for C8690.
$hyperlink = $cell.Hyperlinks.Item(1)
$hyperlink.Address = "https://lcsc.com/product-detail/USB_FTDI_FT232RL-REEL_FT232RL-REEL_C8690.html/?href=jlc-SMT"

# Reflect the edit in the UI selection state (last cell touched/selected).
$cell.Select()
